$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's results as a new row (daily update).
$row = 63
$ws.Cells.Item($row, 1).Value = 46012
$ws.Cells.Item($row, 2).Value = 134
$ws.Cells.Item($row, 3).Value = 150
$ws.Cells.Item($row, 4).Value = 138

# Match the date-style formatting used by the other rows in column A.
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
